# UNO seats per room.xlsx - reformat sheet for VLOOKUP usage.
# Inserts a helper "Room Helper" column between Room (A) and Seats (B/now C)
# that mirrors the Room value as text via TEXT(), so a VLOOKUP against a
# text-typed key (e.g. from the crosslisted data) will match correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; this shifts the old "Seats" column from B to C.
$ws.Columns("B:B").Insert()

# New column header.
$ws.Range("B1").Value = "Room Helper"

# B2 gets its own (unshared) formula, matching how a user would type the
# first formula directly into the cell.
$ws.Range("B2").Formula = "=TEXT(A2,0)"

# B3:B42 are filled from B2 as a single shared-formula block.
$ws.Range("B3:B42").Formula = "=TEXT(A3,0)"

# Row 42's Room value (A42) becomes a formula returning "248" as well,
# consistent with the rest of the sheet being driven off of formulas.
$ws.Range("A42").Formula = "=TEXT(248,0)"

# Update the visible selection to match the new helper column.
$ws.Range("B2:B42").Select()
